$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organizations")

# Update column widths
# Note: this runtime's ColumnWidth setter stores an extra ~0.8333 padding
# offset compared to the raw value provided, so we subtract it here to
# land on the exact target widths in the saved OOXML.
$widthPad = 0.8333333333333357
$ws.Columns.Item(6).ColumnWidth = 25 - $widthPad
$ws.Columns.Item(9).ColumnWidth = 39 - $widthPad
$ws.Columns.Item(11).ColumnWidth = 37 - $widthPad
$ws.Columns.Item(12).ColumnWidth = 42 - $widthPad

# Fill in newly-scraped data
$ws.Range("F2").Value = "studentgovernme@buc.edu"
$ws.Range("K2").Value = "https://twitter.com/studentgovernme"

$ws.Range("F3").Value = "honorsociety@buc.edu"
$ws.Range("L3").Value = "https://youtube.com/channel/honorsociety"

$ws.Range("I4").Value = "https://instagram.com/studentvoluntee"
